# Performance.xlsx update: add cudaGammaCorrection / cudaAntisotropy API rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for two new API blocks (4 new rows) right before the old
#    "Updating" footer rows (old rows 30:31 -> become rows 34:35).
# ---------------------------------------------------------------------------
$ws.Range("A30:E33").Insert(-4121)

# ---------------------------------------------------------------------------
# 2. New row pair for "cudaGammaCorrection" (rows 30-31), styled like the
#    other normal API rows above it (A/E centre-aligned style "2" already
#    carried over by the row-insert above).
# ---------------------------------------------------------------------------
$ws.Range("A30").Value = "cudaGammaCorrection"
$ws.Range("B30").Value = "1024 * 1024"
$ws.Range("D30").Value = 0.876
$ws.Range("E30").Value = "3-ch image"

$ws.Range("B31").Value = "2048 * 2048"
$ws.Range("D31").Value = 3.492

$ws.Range("A30:A31").Merge()
$ws.Range("E30:E31").Merge()

# ---------------------------------------------------------------------------
# 3. New row pair for "cudaAntisotropy" (rows 32-33). These use a distinct,
#    smaller strike-through font (9pt) to flag the entry. Build the format
#    once on A32 / B32 and then fan it out via copy/paste-format so every
#    cell lands on the same two style records instead of minting new ones.
# ---------------------------------------------------------------------------
$ws.Range("A32").Value = "cudaAntisotropy"
$ws.Range("A32").HorizontalAlignment = -4108
$ws.Range("A32").VerticalAlignment = -4108
$ws.Range("A32").Font.Strikethrough = $true
$ws.Range("A32").Font.Size = 9

$ws.Range("B32").Value = "1024 * 1024"
$ws.Range("B32").Font.Strikethrough = $true
$ws.Range("B32").Font.Size = 9

$ws.Range("D32").Value = 16.277
$ws.Range("B33").Value = "2048 * 2048"
$ws.Range("D33").Value = 65.289
$ws.Range("E32").Value = "3-ch image, loop times is 10"

# Fan out the "centre, strike 9pt" format (A32) to A33, E32, E33.
$ws.Range("A32").Copy()
$ws.Range("A33,E32,E33").PasteSpecial(-4122)

# Fan out the "plain, strike 9pt" format (B32) to B33, C32:C33, D32:D33.
$ws.Range("B32").Copy()
$ws.Range("B33,C32,C33,D32,D33").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("A32:A33").Merge()
$ws.Range("E32:E33").Merge()

# ---------------------------------------------------------------------------
# 4. Update the view so the newly added rows are front and centre.
# ---------------------------------------------------------------------------
$ws.Range("A32:E35").Select()

# ---------------------------------------------------------------------------
# 5. Page is set up for portrait printing.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
